# product interface document update
#
# The only functional data change in this revision is a password rotation
# on the "40.73.1.192" sheet: the SSH password for the `neuroot` account
# (row 3, column C "密码") is updated to match the value already used by
# `root` on the row below, i.e. it becomes "qeWf25?Bo".
#
# The revision also reflects that the author left the workbook with the
# first worksheet ("40.73.1.192") active/selected (rather than the second,
# "40.73.0.185", which was active before), so we activate that sheet and
# move its selection accordingly.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("40.73.1.192")

# --- Password rotation: neuroot's ssh password now matches root's ---
$ws1.Range("C3").Value = "qeWf25?Bo"

# --- Make the first sheet the active/selected tab, as in the saved file ---
$ws1.Activate()
$ws1.Range("C18").Select()
